$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.647218666666667
$ws.Range("H2").Value = 7.941656
$ws.Range("J2").Value = 0.4640059894538356
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.523753
$ws.Range("N2").Value = 58.571259
$ws.Range("O2").Value = 0.4652898160202426
$ws.Range("P2").Value = 0.4652898160202426
$ws.Range("Q2").Value = 51.68364338498933
$ws.Range("R2").Value = 465.152790464904
$ws.Range("S2").Value = 0.2158972614652658
$ws.Range("T2").Value = 0.2158972614652658

$ws.Range("G3").Value = 2.647218666666667
$ws.Range("H3").Value = 7.941656
$ws.Range("J3").Value = 0.4640059894538356
$ws.Range("O3").Value = 0.4569298967820781
$ws.Range("P3").Value = 0.4569298967820781
$ws.Range("S3").Value = 0.2120182088674072
$ws.Range("T3").Value = 0.2120182088674072

$ws.Range("G4").Value = 2.647218666666667
$ws.Range("H4").Value = 7.941656
$ws.Range("J4").Value = 0.4640059894538356
$ws.Range("O4").Value = 0.07778028719767933
$ws.Range("P4").Value = 0.07778028719767933
$ws.Range("Q4").Value = 8.639709031869334
$ws.Range("R4").Value = 77.757381286824
$ws.Range("S4").Value = 0.03609051912116271
$ws.Range("T4").Value = 0.03609051912116271

$ws.Range("I5").Value = 0.5359940105461642
$ws.Range("J5").Value = 0.5359940105461642
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.523753
$ws.Range("N5").Value = 58.571259
$ws.Range("O5").Value = 0.4652898160202426
$ws.Range("P5").Value = 0.4652898160202426
$ws.Range("Q5").Value = 59.70208128167766
$ws.Range("R5").Value = 537.3187315350989
$ws.Range("S5").Value = 0.2493925545549767
$ws.Range("T5").Value = 0.2493925545549767

$ws.Range("I6").Value = 0.5359940105461642
$ws.Range("J6").Value = 0.5359940105461642
$ws.Range("O6").Value = 0.4569298967820781
$ws.Range("P6").Value = 0.4569298967820781
$ws.Range("S6").Value = 0.2449116879146709
$ws.Range("T6").Value = 0.2449116879146709

$ws.Range("I7").Value = 0.5359940105461642
$ws.Range("J7").Value = 0.5359940105461642
$ws.Range("O7").Value = 0.07778028719767933
$ws.Range("P7").Value = 0.07778028719767933
$ws.Range("S7").Value = 0.04168976807651662
$ws.Range("T7").Value = 0.04168976807651662
